$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in the new "Obstacles a l'Ecoulement" (ROE) monitoring-programme row (row 20,
# "obstacles_ecoulement"), which was previously present only as a placeholder ("continuites")
# with an empty row.
$ws.Range("A20").Value = "obstacles_ecoulement"
$ws.Range("B20").Value = "Obstacles à l'Ecoulement"
$ws.Range("D20").Value = "Le ROE permet d’avoir une information spatialisée sur les obstacles à l’écoulement des cours d’eau d’origine humaine (caractéristiques, usages, gestion). Cette information peut être complétée dans la BDOE."
$ws.Range("E20").Value = "Collecter des informations objectives sur le nombre, la localisation et les caractéristiques des obstacles à l'écoulement sur les cours d'eau"
$ws.Range("F20").Value = "Elaboration d’un référentiel national et un socle commun d’information. Estimer la pression ouvrage sur les cours d'eau. `nSuivi des politiques de restauration de la continuité écologique`nCalcul d’indicateurs de continuité écologique"
$ws.Range("G20").Value = "75,77,78,91,92,93,94,95"
$ws.Range("J20").Value = "Inventaire:1,2,3,4,5,6,7,8,9,10,11,12`nHauteur de chute:5,6,7,8,9"
$ws.Range("K20").Value = "Opportuniste ou plannifié`nA l'étiage pour la caractérisation des hauteurs de chute"
$ws.Range("L20").Value = "Animation nationale: Karl Kreutzenberger`nAnimation régionale:`nCédric Mondy`nAssistance:`nassistance.geobs@ofb.gouv.fr"
$ws.Range("M20").Value = "Agence de l'eau`nDRIEAT`nSyndicats de rivière`nFédérations de pêche"
$ws.Range("O20").Value = "variable"
$ws.Range("P20").Value = "Pratique"
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Animation`nBase de données"
$ws.Range("S20").Value = "Coordination`nContribution à l'évaluation des politiques publiques mises en place au niveau du bassin`nFormation"
$ws.Range("T20").Value = "Animation locale`nSaisie`nValidation"
$ws.Range("U20").Value = "Outil Géobs:`n- Consultation: sur demande`n- Saisie: Information dispensée par les agents de la DR`n- Validation: Correspondants continuité formés"
$ws.Range("W20").Value = "'- Fiche terrain`n- GPS`n- Appareil photo`n- Mire ou autre équipement permettant de mesurer la hauteur de chute"
$ws.Range("Y20").Value = "Saisie sur GéObs (https://geobs.eaufrance.fr/)`nValidation obligatoire par un agent de l'OFB disposant du profil validation"
$ws.Range("Z20").Value = "Sandre:`nhttp://www.sandre.eaufrance.fr"
$ws.Range("AA20").Value = "texte:La continuité écologique des cours d'eau;lien:https://www.ofb.gouv.fr/la-continuite-ecologique-des-cours-deau"
$ws.Range("AB20").Value = "texte:Dataviz nationale;lien:https://professionnels.ofb.fr/fr/doc-dataviz/dataviz-mieux-connaitre-ouvrages-qui-jalonnent-nos-cours-deau"
$ws.Range("AD20").Value = "texte:Bilans (serveur DR);lien:\\ad.intra\dfs\COMMUNS\REGIONS\IDF\DR\05_CONNAISSANCE\ROE\04_Bilans"
$ws.Range("AE20").Value = "texte:La méthode ICE;lien:https://professionnels.ofb.fr/fr/node/387"

# V20 was left with a stale "hyperlink" cell style (blue/underlined) even though it is
# empty; clear that formatting back to the plain wrapped-text style used by its neighbours.
$ws.Range("V20").ClearFormats()
$ws.Range("V20").WrapText = $true

# Row 20 now holds a full record (like the other populated rows), so give it the same
# kind of explicit row height used for the other long/wrapped rows.
$ws.Rows.Item(20).RowHeight = 165

# Restore the view to a plain, non-edited state (matches the final author selection).
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("X20").Select()
